$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.586.15'
$ws.Range("E2").Value = '  +0.58%  '
$ws.Range("D3").Value = '1.822.47'
$ws.Range("E3").Value = '  +1.48%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.008'
$ws.Range("E5").Value = '  +0.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '305.24'
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4645'
$ws.Range("E7").Value = '  +2.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3594'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07131'
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8981'
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07770'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.32'
$ws.Range("E12").Value = '  -1.09%  '
$ws.Range("D13").Value = '1.795.77'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.248'
$ws.Range("E14").Value = '  -0.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.311'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.20'
$ws.Range("E16").Value = '  +2.48%  '
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008536'
$ws.Range("E18").Value = '  -0.61%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '26.623.13'
$ws.Range("E20").Value = '  +0.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.11'
$ws.Range("E21").Value = '  -1.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.011'
$ws.Range("E22").Value = '  +0.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.53'
$ws.Range("E23").Value = '  -0.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.915'
$ws.Range("E24").Value = '  -3.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.95'
$ws.Range("E25").Value = '  -0.50%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.88'
$ws.Range("E26").Value = '  -0.28%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.967'
$ws.Range("E27").Value = '  -4.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '113.63'
$ws.Range("E28").Value = '  +1.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.797'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08796'
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.141'
$ws.Range("E31").Value = '  +2.64%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7291'
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.722'
$ws.Range("E33").Value = '  -1.37%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.430'
$ws.Range("E34").Value = '  -0.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.127'
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.073'
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01920'
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.918'
$ws.Range("E38").Value = '  +1.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05095'
$ws.Range("E39").Value = '  -0.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.850'
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5027'
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1490'
$ws.Range("E42").Value = '  -1.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.942'
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.009'
$ws.Range("E44").Value = '  +0.10%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4637'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.972'
$ws.Range("E46").Value = '  +0.63%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '98.15'
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.553'
$ws.Range("E48").Value = '  -2.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05979'
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.60'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.71'
$ws.Range("E51").Value = '  -1.69%  '
